# Add "ankita" login credentials as a new worksheet "TestData-Ankita"
$wb = $excel.ActiveWorkbook

# Leave the "TestData-Rutuja" sheet's selection where the author last left it
# before adding the new sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Select() | Out-Null
$lastSheet.Range("F11").Select() | Out-Null

# Add the new worksheet after the last existing sheet ("TestData-Rutuja")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "TestData-Ankita"

# Populate cells in the same order the original author typed them, so the
# shared-string table is built up in the same sequence as the source file.
$newSheet.Range("A1").Value = "Username"
$newSheet.Range("B1").Value = "Passward"
$newSheet.Range("C1").Value = "Invalid_Username"

$newSheet.Range("A2").Value = "ankita"
$newSheet.Range("B2").Value = "Pass9Test"
$newSheet.Range("C2").Value = "ankita987"
$newSheet.Range("D2").Value = "ankita987"

$newSheet.Range("D1").Value = "Invalid_Password"

# Make the new sheet the active/selected tab with D1 selected
$newSheet.Select() | Out-Null
$newSheet.Range("D1").Select() | Out-Null
